$d = $word.ActiveDocument

$replacements = @(
    @("2024-01-11 Thursday", "2024-01-12 Friday"),
    @("884×8=7072", "374×5=1870"),
    @("313×5=1565", "996×9=8964"),
    @("120×5=600", "277×4=1108"),
    @("556×6=3336", "177×3=531"),
    @("709×9=6381", "409×2=818"),
    @("288×6=1728", "117×2=234"),
    @("154×6=924", "870×7=6090"),
    @("358×4=1432", "442×2=884"),
    @("397×7=2779", "632×5=3160"),
    @("121×3=363", "620×4=2480"),
    @("949×6=5694", "666×2=1332"),
    @("997×4=3988", "809×7=5663"),
    @("459×7=3213", "475×4=1900"),
    @("493×6=2958", "454×3=1362"),
    @("192×8=1536", "985×9=8865"),
    @("552×6=3312", "174×9=1566"),
    @("101×2=202", "205×9=1845"),
    @("848×6=5088", "197×3=591"),
    @("857×2=1714", "477×2=954"),
    @("287×3=861", "443×9=3987"),
    @("582×2=1164", "390×5=1950"),
    @("238×9=2142", "750×6=4500"),
    @("402×3=1206", "660×3=1980"),
    @("874×4=3496", "718×2=1436"),
    @("261×4=1044", "843×2=1686")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
